# Phase 4 rows (63-79) of the WBS sheet: mark as Done, stamp the completion
# date, and fill in the evidence checklist columns (Schema / Validation /
# Permissions-Isolation / Workflow / Evidence).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("WBS")

$check = [char]0x2705   # ✅

for ($row = 63; $row -le 79; $row++) {
    $ws.Cells.Item($row, 8).Value = "Done"    # H: Execution Status

    # K: Completed On -- copy J (Started On), which already holds the
    # "2026-05-01" text value with the correct date-style formatting, so K
    # picks up the same style/text representation instead of Excel
    # auto-converting the string into a date serial number.
    $ws.Cells.Item($row, 10).Copy($ws.Cells.Item($row, 11))

    $ws.Cells.Item($row, 12).Value = $check   # L: Schema
    $ws.Cells.Item($row, 13).Value = $check   # M: Validation
    $ws.Cells.Item($row, 14).Value = $check   # N: Permissions/Isolation
    $ws.Cells.Item($row, 15).Value = $check   # O: Workflow
    $ws.Cells.Item($row, 16).Value = $check   # P: Evidence
}
